# The presentation's slide-master theme ("Integral" / "Red Violet" colour
# scheme, stored in ppt/theme/theme1.xml) is switched to the default
# "Office Theme" colour scheme (the palette that previously only lived in
# ppt/theme/theme2.xml, used by the Notes Master).
#
# PowerPoint's ColorScheme.Colors(index).RGB slots map onto the 12 theme
# colours in this fixed order:
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink
#
# RGB() in the PowerPoint/VBA object model packs a colour as
# 0x00BBGGRR (blue in the high byte), so convert each target hex colour
# accordingly before assigning it.

function ConvertTo-OleRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $scheme.Colors($i).RGB = ConvertTo-OleRgb $officeThemeColors[$i - 1]
}
